# Apply the LOT2046 course-sheet update:
#  - Ativação date: 01/01/2018 -> 01/01/2022
#  - Add missing English "Objectives:" body text (row 11, columns B & C)
#  - Rewrite "Programa resumido:" (row 14) to drop the cellular-structure clause
#  - Rewrite "Short syllabus:" (row 15) to drop the cellular-structure clause
#  - Rewrite "Programa:" (row 16) to replace the opening clause with the new
#    "Estrutura e função das principais moléculas orgânicas..." clause
#  - Rewrite "Syllabus:" (row 17) with the matching English opening clause

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2018 -> 01/01/2022
# (leading apostrophe forces text so Excel doesn't silently convert the
#  date-shaped string into a date serial number)
$ws.Range("B8").Value2 = "'01/01/2022"
$ws.Range("C8").Value2 = "'01/01/2022"

# Row 11 ("Objectives:") gained the English translation of row 10's body,
# mirroring the B/C layout already used by every other label row (B = plain
# "current" text, C = "modified" text shown in red).
$objectivesEn = "Provide the necessary knowledge on the fundamental aspects of Microbiology and Microbial Biochemistry and its importance in studies on Ecology of Microorganisms. Provide knowledge about the role and use of microorganisms in biological processes of interest to Environmental Engineering."
$ws.Range("B11").Value2 = $objectivesEn
$ws.Range("C11").Value2 = $objectivesEn

# New cells inherit the row's style (bold, from column A's label cell) by
# default; restore the normal look used by the rest of column B/C.
$ws.Range("B11").WrapText = $true
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").Font.Bold = $false
$ws.Range("C11").WrapText = $true
$ws.Range("C11").VerticalAlignment = -4160
$ws.Range("C11").Font.Color = 255

# Programa resumido: drop the "Estrutura celular e história evolutiva;" clause
$programaResumido = "Diversidade metabólica; cultivo e crescimento microbiano; isolamento microbiano; ecossistemas microbianos; biorremediação e biodeterioração  microbiana; bioindicadores."
$ws.Range("B14").Value2 = $programaResumido
$ws.Range("C14").Value2 = $programaResumido

# Short syllabus: drop the "Cellular structure and evolutive history;" clause
$shortSyllabus = "Metabolic diversity; microbial culture and growth; microbial isolation; microbial ecosystems; microbial bioremediation and biodeterioration; bioindicators."
$ws.Range("B15").Value2 = $shortSyllabus
$ws.Range("C15").Value2 = $shortSyllabus

# Programa: replace the opening clause about cellular structure with the new
# clause about the structure/function of the main organic molecules.
$programa = "Estrutura e função das principais moléculas orgânicas: carboidratos, lipídeos,proteínas e ácidos nucleicos.–Diversidade metabólica: Micro-organismos autotróficos e heterotróficos; glicólise; fermentações; respiração; via das pentoses-fosfato; fotossíntese. –Cultivo e crescimento microbiano: Nutrição microbiana; meios de cultura; fatores ambientais; reprodução e crescimento; medidas e controle de crescimento microbiano. –Isolamento microbiano: Técnicas e meios de isolamento.–Ecossistemas microbianos: Diversidade microbiana e ciclos biogeoquímicos. –Biorremediação e biodeterioração microbiana: Lixiviação bacteriana de metais; bioacumulação e biotransformação microbiana de metais; biodegradação de materiais lignocelulósicos; biodegradação de hidrocarbonetos; biodeterioração de monumentos históricos. –Bioindicadores: Bioindicadores de qualidade de água, ar e solo."
$ws.Range("B16").Value2 = $programa
$ws.Range("C16").Value2 = $programa

# Syllabus: same substitution in English.
$syllabus = "Structure and function of the main organic molecules: carbohydrates, lipids,proteins and nucleic acids.Metabolic diversity: autotrophic and heterotrophic microorganisms, glycolysis; fermentations; respiration; pentose-phosphate pathway; photosynthesis. Microbial culture and growth: microbial nutrition; culture media; ambiental factors; reproduction and growth; measures and control of microbial growth.Microbial isolation: techniques and isolation media.Microbial ecosystems: microbial diversity and biogeochemical cycles.  Microbial biorremediation and biorremediation: bacterial leaching of metals; microbial bioaccumulation and biotransformation of metals; biodegradation of lignocellulosic materials; biodegradation of hydrocarbonets; biodeterioration of hystoric monuments. Bioindicators: bioindicators of the quality of water, air and soil."
$ws.Range("B17").Value2 = $syllabus
$ws.Range("C17").Value2 = $syllabus
